$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values are prefixed with a leading apostrophe (PowerShell single-quoted
# string starting with two single quotes -> one literal apostrophe) so that
# Excel stores them as literal text instead of auto-converting number-like
# strings (e.g. "43.30", "0.679") into floating point numbers, which would
# drop trailing zeros / alter precision and change the cell type.
$ws.Range('D2').Value = '''35.238.39'
$ws.Range('E2').Value = '''  +0.32%  '
$ws.Range('D3').Value = '''1.876.85'
$ws.Range('E3').Value = '''  -1.41%  '
$ws.Range('E4').Value = '''  -0.54%  '
$ws.Range('D5').Value = '''245.41'
$ws.Range('E5').Value = '''  -3.28%  '
$ws.Range('D6').Value = '''0.679'
$ws.Range('E6').Value = '''  -3.25%  '
$ws.Range('E7').Value = '''  -0.61%  '
$ws.Range('D8').Value = '''43.30'
$ws.Range('E8').Value = '''  +4.35%  '
$ws.Range('D9').Value = '''0.353'
$ws.Range('E9').Value = '''  -1.68%  '
$ws.Range('D10').Value = '''53.52'
$ws.Range('E10').Value = '''  +2.09%  '
$ws.Range('D11').Value = '''0.0734'
$ws.Range('E11').Value = '''  -2.91%  '
$ws.Range('E12').Value = '''  -0.57%  '
$ws.Range('D13').Value = '''13.41'
$ws.Range('E13').Value = '''  +1.77%  '
$ws.Range('D14').Value = '''2.149.93'
$ws.Range('E14').Value = '''  -1.41%  '
$ws.Range('D15').Value = '''0.766'
$ws.Range('E15').Value = '''  +4.67%  '
$ws.Range('B16').Value = '''WrappedEther'
$ws.Range('C16').Value = '''https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '''1.901.49'
$ws.Range('E16').Value = '''  -0.11%  '
$ws.Range('B17').Value = '''Polkadot'
$ws.Range('C17').Value = '''https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D17').Value = '''4.88'
$ws.Range('E17').Value = '''  -2.89%  '
$ws.Range('D18').Value = '''35.246.59'
$ws.Range('E18').Value = '''  +0.33%  '
$ws.Range('D19').Value = '''72.64'
$ws.Range('E19').Value = '''  -1.58%  '
$ws.Range('D20').Value = '''0.0₃0817'
$ws.Range('E20').Value = '''  -2.90%  '
$ws.Range('D21').Value = '''243.03'
$ws.Range('E21').Value = '''  +0.09%  '
$ws.Range('D22').Value = '''12.72'
$ws.Range('E22').Value = '''  -2.29%  '
$ws.Range('D23').Value = '''4.93'
$ws.Range('E23').Value = '''  -2.39%  '
$ws.Range('D24').Value = '''2.63'
$ws.Range('E24').Value = '''  +7.48%  '
$ws.Range('E25').Value = '''  -0.64%  '
$ws.Range('E26').Value = '''  -7.43%  '
$ws.Range('D27').Value = '''165.45'
$ws.Range('E27').Value = '''  -1.11%  '
$ws.Range('D28').Value = '''8.47'
$ws.Range('E28').Value = '''  -1.54%  '
$ws.Range('D29').Value = '''18.18'
$ws.Range('E29').Value = '''  -1.97%  '
$ws.Range('D30').Value = '''0.126'
$ws.Range('E30').Value = '''  -2.86%  '
$ws.Range('D31').Value = '''4.128.45'
$ws.Range('E31').Value = '''  -0.01%  '
$ws.Range('E32').Value = '''  +6.59%  '
$ws.Range('D33').Value = '''2.00'
$ws.Range('E33').Value = '''  -0.70%  '
$ws.Range('E34').Value = '''  -1.80%  '
$ws.Range('D35').Value = '''0.0588'
$ws.Range('E35').Value = '''  -2.55%  '
$ws.Range('E36').Value = '''  -2.70%  '
$ws.Range('E37').Value = '''  -0.58%  '
$ws.Range('D38').Value = '''0.834'
$ws.Range('E38').Value = '''  -2.17%  '
$ws.Range('D39').Value = '''0.0733'
$ws.Range('E39').Value = '''  +12.48%  '
$ws.Range('E40').Value = '''  -4.25%  '
$ws.Range('D41').Value = '''17.61'
$ws.Range('E41').Value = '''  +2.06%  '
$ws.Range('E42').Value = '''  +0.16%  '
$ws.Range('D43').Value = '''95.76'
$ws.Range('E43').Value = '''  -4.76%  '
$ws.Range('E44').Value = '''  -3.35%  '
$ws.Range('D45').Value = '''1.301.81'
$ws.Range('E45').Value = '''  -0.10%  '
$ws.Range('E46').Value = '''  -2.52%  '
$ws.Range('D47').Value = '''0.0794'
$ws.Range('E47').Value = '''  +5.59%  '
$ws.Range('E48').Value = '''  -1.93%  '
$ws.Range('E49').Value = '''  -1.07%  '
$ws.Range('D50').Value = '''11.94'
$ws.Range('E50').Value = '''  -4.28%  '
$ws.Range('D51').Value = '''6.20'
$ws.Range('E51').Value = '''  -5.89%  '
